$wb = $excel.ActiveWorkbook

# Add new worksheet "tc002" after the last existing sheet
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "tc002"

# Header row
$ws.Range("A1").Value = "fieldName"
$ws.Range("B1").Value = "textBox"
$ws.Range("C1").Value = "objType1"
$ws.Range("D1").Value = "objType2"

# Data row
$ws.Range("A2").Value = "Req"
$ws.Range("B2").Value = "Text Box"
$ws.Range("C2").Value = "Requirements"
$ws.Range("D2").Value = "Test Cases"

# Column widths (approximate engine's character-width -> stored-width rounding)
$ws.Columns.Item(1).ColumnWidth = 12.1
$ws.Columns.Item(2).ColumnWidth = 11.9
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(4).ColumnWidth = 14.75

# Row heights
$ws.Rows.Item(1).RowHeight = 20
$ws.Rows.Item(2).RowHeight = 20.5

# Apply alignment (vertical centered + wrap text) to the data range as a single
# combined style, by building the style on the first cell and then copying the
# format across the rest of the range (avoids leaving stray intermediate styles
# behind when setting two alignment properties cell-by-cell).
$formatCell = $ws.Range("A1")
$formatCell.WrapText = $true
$formatCell.VerticalAlignment = -4108
$formatCell.Copy()
$ws.Range("A1:D2").PasteSpecial(-4122)

# Selection / active cell for the new sheet
$ws.Activate()
$ws.Range("D6").Select()

# Active tab becomes the new sheet (index 2, 0-based) once selected/activated
Write-Output "done"
